$wb = $excel.ActiveWorkbook

# The same data updates need to be applied to both the "展览" and
# "全部类型" sheets (they mirror the same underlying data).
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F3").Value  = 3100
    $ws.Range("F5").Value  = 2677
    $ws.Range("F9").Value  = 1415
    $ws.Range("F13").Value = 1206
    $ws.Range("F14").Value = 5
    $ws.Range("F15").Value = 363
    $ws.Range("F17").Value = 39
    $ws.Range("F18").Value = 36
    $ws.Range("G18").Value = 50
    $ws.Range("F22").Value = 2590
    $ws.Range("F24").Value = 296
}
